# Apply a cyclic permutation of the data rows (2-10) on the active sheet.
# The record previously on a given row moves to a new row as follows
# (derived from the target OOXML diff):
#   old row 2  -> new row 8
#   old row 3  -> new row 2
#   old row 4  -> new row 3
#   old row 5  -> new row 4
#   old row 8  -> new row 9
#   old row 9  -> new row 10
#   old row 10 -> new row 5
#   old row 6  -> new row 6   (unchanged)
#   old row 7  -> new row 7   (unchanged)
#
# Implementation notes:
#  - Cell.Copy is used instead of re-assigning .Value/.Value2, so that
#    cell typing is preserved exactly as authored. In particular the
#    Start/Slut date & time columns are stored as plain text
#    ("2023-03-23"); assigning such a string back through .Value/.Value2
#    makes Excel "helpfully" reinterpret it as a real date serial
#    number, which does not match the source workbook.
#  - Work is done cell-by-cell (not whole-row Range.Copy) because
#    copying a whole A:AY row range pads every column of the
#    destination with an explicit (blank) cell, even for columns that
#    had no <c> element at all in the original sheet.
#  - Because this is a permutation (several rows trade places at once)
#    each moving row is first copied to a scratch area below the used
#    range; the scratch copies are then copied into their final
#    destination. Whether a destination cell should hold a copied value
#    is decided from a snapshot of "does the original source cell have
#    a value" taken up front, before anything is modified -- re-reading
#    .Value2 after a cell has already been copied once is not reliable
#    for this in this engine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1    # A
$lastCol  = 51   # AY

# Rows that move (row 6 and row 7 keep their data in place).
$movingRows = @(2, 3, 4, 5, 8, 9, 10)

# Destination row -> source row, for the moving rows.
$destToSrc = @{
    2  = 3
    3  = 4
    4  = 5
    5  = 10
    8  = 2
    9  = 8
    10 = 9
}

# Scratch rows (well below the used range A1:AY10) used to stage a copy
# of each moving row before any original cell gets overwritten/cleared.
$scratchRow = @{
    2  = 101
    3  = 102
    4  = 103
    5  = 104
    8  = 105
    9  = 106
    10 = 107
}

# 0. Snapshot, for every moving row/column, whether the *original* cell
#    holds anything at all (including an "empty" inline string). This
#    is taken before any writes happen, and is what later steps trust.
$hasValue = @{}
foreach ($row in $movingRows) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $hasValue[[string]$row + "_" + [string]$col] = ($cell.Value2 -ne $null)
    }
}

# 1. Stage a copy of every moving row, cell by cell, only where the
#    original had content.
foreach ($row in $movingRows) {
    $scratch = $scratchRow[$row]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        if ($hasValue[[string]$row + "_" + [string]$col]) {
            $srcCell = $ws.Cells.Item($row, $col)
            $dstCell = $ws.Cells.Item($scratch, $col)
            $srcCell.Copy($dstCell)
        }
    }
}

# 2. Clear the original rows so that columns which must end up blank
#    (because the mapped source had no value there) really are blank --
#    pasting/copying "nothing" onto a non-blank cell is a no-op rather
#    than a clear in this engine.
foreach ($row in $movingRows) {
    $rowRange = $ws.Range($ws.Cells.Item($row, $firstCol), $ws.Cells.Item($row, $lastCol))
    $rowRange.Clear()
}

# 3. Write the staged data into its final destination row, using the
#    up-front snapshot to know which cells should actually receive a
#    value.
foreach ($destRow in $destToSrc.Keys) {
    $srcRow = $destToSrc[$destRow]
    $scratch = $scratchRow[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        if ($hasValue[[string]$srcRow + "_" + [string]$col]) {
            $stagedCell = $ws.Cells.Item($scratch, $col)
            $destCell = $ws.Cells.Item($destRow, $col)
            $stagedCell.Copy($destCell)
        }
    }
}

# 4. Clean up the scratch area.
foreach ($row in $movingRows) {
    $scratch = $scratchRow[$row]
    $scratchRange = $ws.Range($ws.Cells.Item($scratch, $firstCol), $ws.Cells.Item($scratch, $lastCol))
    $scratchRange.Clear()
}
